# FINFLUX-2789 Correcting Failed Overdue scenarios
#
# - Adds a new "Modify Transaction1" sheet (OverDueTillDate / clickonsubmit /
#   NavigateToLoan helper rows) used to drive the overdue-correction automation.
# - Tweaks the "Transactions" sheet selection/column width.
# - Leaves "Repayment schedule" as no-longer-the-active-tab (the new sheet
#   becomes the active / selected tab instead).

$wb = $excel.ActiveWorkbook

$loanInput = $wb.Worksheets.Item("NewLoanInput")
$transactions = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------------
# 1. Transactions sheet: widen column A slightly and move the selection.
# ---------------------------------------------------------------------------
$transactions.Columns.Item(1).ColumnWidth = 3.2
$transactions.Range("E8").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Modify Transaction1" sheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Modify Transaction1"

# Column A (labels) - reuse the header-style formatting (fill + font) used
# elsewhere in the workbook.
$loanInput.Range("A1").Copy()
$newSheet.Range("A1:A3").PasteSpecial(-4122)
$newSheet.Range("A1:A3").Font.Name = "Calibri"

# Column B, rows 2 & 3 (values) - reuse the flat/date style.
$loanInput.Range("B3").Copy()
$newSheet.Range("B2:B3").PasteSpecial(-4122)
$newSheet.Range("B2:B3").WrapText = $False

# Column B, row 1 (date) - reuse the date-formatted style as-is.
$newSheet.Range("B1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "OverDueTillDate"
$newSheet.Range("B1").Value = 42064
$newSheet.Range("A2").Value = "clickonsubmit"
$newSheet.Range("B2").Value = "Submit"
$newSheet.Range("A3").Value = "NavigateToLoan"
$newSheet.Range("B3").Value = "navigate"

$newSheet.Columns.Item(1).ColumnWidth = 15
$newSheet.Columns.Item(2).ColumnWidth = 12.6

$newSheet.Range("A1:B3").Select()
